$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "301.42"
Set-TextValue "E2" "0.90%"
Set-TextValue "D3" "31.65"
Set-TextValue "E3" "1.45%"
Set-TextValue "D4" "5.089"
Set-TextValue "E4" "-1.26%"
Set-TextValue "D5" "0.07793"
Set-TextValue "E5" "-2.98%"
Set-TextValue "D6" "2.232"
Set-TextValue "E6" "-17.95%"
Set-TextValue "D7" "7.794"
Set-TextValue "E7" "-0.39%"
Set-TextValue "D8" "3.816"
Set-TextValue "E8" "-0.20%"
Set-TextValue "D9" "0.9180"
Set-TextValue "E9" "0.27%"
Set-TextValue "D10" "0.1754"
Set-TextValue "E10" "1.00%"
Set-TextValue "D11" "0.07548"
Set-TextValue "D12" "0.08993"
Set-TextValue "E12" "8.07%"
Set-TextValue "D13" "0.03026"
Set-TextValue "E13" "1.04%"
Set-TextValue "D14" "0.1002"
Set-TextValue "E14" "0.67%"
Set-TextValue "E15" "1.27%"
Set-TextValue "D16" "0.006055"
Set-TextValue "E16" "-0.89%"
Set-TextValue "D17" "3.468"
Set-TextValue "E17" "-0.89%"
Set-TextValue "E18" "-0.01%"
Set-TextValue "D19" "0.3292"
Set-TextValue "E19" "0.26%"
Set-TextValue "E20" "0.67%"
Set-TextValue "E21" "-6.48%"
Set-TextValue "D23" "0.04594"
Set-TextValue "E23" "0.35%"
Set-TextValue "D24" "0.001251"
Set-TextValue "E24" "-0.68%"
Set-TextValue "D25" "0.004472"
Set-TextValue "E25" "0.64%"
Set-TextValue "E26" "5.72%"
Set-TextValue "E27" "-1.49%"
Set-TextValue "E39" "-3.19%"
Set-TextValue "D40" "0.04780"
Set-TextValue "E40" "5.91%"
Set-TextValue "D41" "0.007528"
Set-TextValue "E41" "7.25%"
Set-TextValue "D42" "0.1359"
Set-TextValue "E42" "1.19%"
Set-TextValue "E43" "-2.43%"
Set-TextValue "D44" "0.01026"
Set-TextValue "E44" "4.34%"
Set-TextValue "D45" "0.00006210"
Set-TextValue "E45" "-4.08%"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "E46" "-0.20%"
Set-TextValue "E47" "28.70%"
Set-TextValue "D48" "0.7335"
Set-TextValue "E48" "-10.61%"
Set-TextValue "D49" "0.00002098"
Set-TextValue "E49" "-0.20%"
Set-TextValue "D50" "0.0001998"
Set-TextValue "E50" "-0.20%"

Write-Host "Updated symbol list values."
